{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// Two changes, matching the target diff:\n//  1) The \"NOTAS IMPORTANTES:\" heading paragraph becomes bold, red (FF0000)\n//     and size 36 half-points (18pt), up from size 24 half-points (12pt).\n//  2) The last paragraph (\"El ejemplo de validaci\u00f3n...\") had its text split\n//     across two runs around a bookmark (_GoBack) purely as an artifact of\n//     editing; the run split is removed so the whole sentence lives in a\n//     single run immediately before the (unchanged) bookmark.\n\nconst body = context.document.body;\n\n// --- Change 1: \"NOTAS IMPORTANTES:\" heading formatting -------------------\nconst notasResults = body.search(\"NOTAS IMPORTANTES:\", { matchCase: true });\nnotasResults.load(\"items\");\nawait context.sync();\n\nconst notasParagraph = notasResults.items[0].paragraphs.getFirst();\nnotasParagraph.font.color = \"#FF0000\";\nnotasParagraph.font.size = 18; // 36 half-points\nawait context.sync();\n\n// --- Change 2: merge the split runs in the final paragraph ---------------\nconst tailText = \"r las validaciones para los dem\u00e1s tipos de campos. \";\nconst headSearch = body.search(\"por defini\", { matchCase: true });\nheadSearch.load(\"items\");\nawait context.sync();\n\n// Move the tail fragment's text so it directly follows \"por defini\",\n// turning the sentence into a single contiguous run ahead of the bookmark.\nconst headRange = headSearch.items[0];\nconst insertionPoint = headRange.getRange(\"End\");\ninsertionPoint.insertText(tailText, \"Start\");\nawait context.sync();\n\n// Two copies of the tail text now exist: the merged one (just inserted)\n// and the original leftover run after the bookmark. Remove the latter.\nconst tailResults = body.search(tailText, { matchCase: true });\ntailResults.load(\"items\");\nawait context.sync();\n\ntailResults.items[tailResults.items.length - 1].delete();\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is already open as $d below.\n#\n# Two changes, matching the target diff:\n#  1) The \"NOTAS IMPORTANTES:\" heading paragraph becomes bold, red (FF0000)\n#     and size 36 half-points (18pt), up from size 24 half-points (12pt).\n#  2) The last paragraph (\"El ejemplo de validaci\u00f3n...\") had its text split\n#     across two runs around a bookmark (_GoBack) purely as an artifact of\n#     editing; the run split is removed so the whole sentence lives in a\n#     single run immediately before the (unchanged) bookmark.\n\n$d = $word.ActiveDocument\n\n# --- Change 1: \"NOTAS IMPORTANTES:\" heading formatting ---------------------\n# Walk the Paragraphs collection (rather than Find, whose collapsed match\n# range does not reliably carry the paragraph-mark's own rPr) so that both\n# the paragraph mark's formatting (w:pPr/w:rPr) and the run's formatting\n# (w:r/w:rPr) get the new color/size, matching Word's own \"select whole\n# paragraph, then format\" behavior.\n$headingParagraph = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $candidate = $d.Paragraphs.Item($i)\n    if ($candidate.Range.Text -like \"*NOTAS IMPORTANTES:*\") {\n        $headingParagraph = $candidate\n        break\n    }\n}\nif ($headingParagraph) {\n    $headingParagraph.Range.Font.Color = 255      # wdColor RGB(255,0,0) -> FF0000\n    $headingParagraph.Range.Font.Size = 18        # 18pt == <w:sz w:val=\"36\"/>\n}\n\n# --- Change 2: merge the split runs in the final paragraph -----------------\n$tailText = \"r las validaciones para los dem\u00e1s tipos de campos. \"\n\n$headRange = $d.Content\n$headRange.Find.ClearFormatting()\n$headFound = $headRange.Find.Execute(\"por defini\")\nif ($headFound) {\n    # Collapse to the end of \"por defini\" and splice the continuation text\n    # in right there, ahead of the bookmark, so the sentence becomes a\n    # single contiguous run.\n    $insertionPoint = $headRange.Duplicate\n    $insertionPoint.Collapse(0)  # wdCollapseEnd\n    $insertionPoint.InsertBefore($tailText)\n\n    # Two copies of the tail text now exist: the merged one (just spliced\n    # in) and the original leftover run after the bookmark. Locate and\n    # delete the latter by searching after the first occurrence.\n    $firstTailRange = $d.Content\n    $firstTailRange.Find.ClearFormatting()\n    $firstTailRange.Find.Execute($tailText) | Out-Null\n\n    $secondTailRange = $d.Content\n    $secondTailRange.Start = $firstTailRange.End\n    $secondTailRange.Find.ClearFormatting()\n    $secondTailFound = $secondTailRange.Find.Execute($tailText)\n    if ($secondTailFound) {\n        $secondTailRange.Delete()\n    }\n}\n"}
